$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 527, shifting rows 527..622 down to 528..623
$ws.Rows.Item(527).Insert()

# Populate new row 527 with the new data
$ws.Cells.Item(527, 1).Value = 7
$ws.Cells.Item(527, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(527, 3).Value = "Ñuble"
$ws.Cells.Item(527, 4).Value = 45209
$ws.Cells.Item(527, 5).Value = 16
$ws.Cells.Item(527, 6).Value = 100112023
$ws.Cells.Item(527, 7).Value = "Brócoli"
$ws.Cells.Item(527, 8).Value = "Sin especificar"
$ws.Cells.Item(527, 9).Value = "Primera"
$ws.Cells.Item(527, 10).Value = 600
$ws.Cells.Item(527, 11).Value = 1000
$ws.Cells.Item(527, 12).Value = 1100
$ws.Cells.Item(527, 13).Value = 1050
$ws.Cells.Item(527, 14).Value = "$/unidad"
$ws.Cells.Item(527, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(527, 16).Value = 1050
$ws.Cells.Item(527, 17).Value = 1
$ws.Cells.Item(527, 18).Value = "Hortaliza"
